$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.165540666666667
$ws.Range("N2").Value = 24.496622
$ws.Range("O2").Value = 0.1715865889461355
$ws.Range("P2").Value = 0.1715865889461355
$ws.Range("Q2").Value = 2.744302125722222
$ws.Range("R2").Value = 24.6987191315
$ws.Range("S2").Value = 0.1715865889461355
$ws.Range("T2").Value = 0.1715865889461355

# Row 3 (only derived specificity columns change)
$ws.Range("O3").Value = 0.4625449807101323
$ws.Range("P3").Value = 0.4625449807101323
$ws.Range("S3").Value = 0.4625449807101323
$ws.Range("T3").Value = 0.4625449807101323

# Row 4
$ws.Range("M4").Value = 13.51552533333333
$ws.Range("N4").Value = 40.546576
$ws.Range("O4").Value = 0.2840084918355372
$ws.Range("P4").Value = 0.2840084918355373
$ws.Range("Q4").Value = 4.542342805777777
$ws.Range("R4").Value = 40.881085252
$ws.Range("S4").Value = 0.2840084918355372
$ws.Range("T4").Value = 0.2840084918355373

# Row 5
$ws.Range("M5").Value = 3.895588
$ws.Range("N5").Value = 11.686764
$ws.Range("O5").Value = 0.08185993850819488
$ws.Range("P5").Value = 0.0818599385081949
$ws.Range("Q5").Value = 1.309242200333333
$ws.Range("R5").Value = 11.783179803
$ws.Range("S5").Value = 0.08185993850819488
$ws.Range("T5").Value = 0.0818599385081949
